# Re-process the dimension/measure metadata with the newly curated dimensions.
# The "balance", "ratios" and "pyg" columns (L, AU, BM) were modelled as
# iaest-dimension:* entries pointing at an xlsx concept mapping; they are now
# curated as iaest-measure:* entries instead (xsd:int typed, "medida" kind,
# no external mapping file). The "municipio-nombre" measure column (N) is
# replaced by a proper sdmx-dimension:refArea dimension (URI-Municipio typed),
# matching the existing refArea columns used for provincia/comarca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column L: "balance" -----------------------------------------------
$ws.Range("L2").Value = "iaest-measure:balance"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
$ws.Range("L5").ClearContents()

# --- Column N: "municipio-nombre" -> refArea dimension ------------------
$ws.Range("N2").Value = "sdmx-dimension:refArea"
$ws.Range("N3").Value = "dim"
$ws.Range("N4").Value = "URI-Municipio"

# --- Column AU: "ratios" -------------------------------------------------
$ws.Range("AU2").Value = "iaest-measure:ratios"
$ws.Range("AU3").Value = "medida"
$ws.Range("AU4").Value = "xsd:int"
$ws.Range("AU5").ClearContents()

# --- Column BM: "pyg" ------------------------------------------------------
$ws.Range("BM2").Value = "iaest-measure:pyg"
$ws.Range("BM3").Value = "medida"
$ws.Range("BM4").Value = "xsd:int"
$ws.Range("BM5").ClearContents()
